$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A4").Value = "2022-11-16_20-48-47"
$ws.Range("B4").Value = " "
$ws.Range("C4").Value = " "
$ws.Range("D4").Value = " "
$ws.Range("E4").Value = "'0.9346306920051575"
$ws.Range("F4").Value = "./checkpoints/2022-11-16_20-48-47"
$ws.Range("G4").Value = "./record/2022-11-16_20-48-47"
$ws.Range("H4").Value = "./log/2022-11-16_20-48-47"
$ws.Range("I4").Value = "./log/2022-11-16_20-48-47/log.txt"

# Clear the implicit "quote prefix" text style Excel applies to E4 so the
# new row's cells stay on the default (unstyled) format, matching rows 2-3.
$ws.Range("A4:I4").Style = "Normal"
